$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their textual representation (matches the
# feed format, e.g. "0.610" or "66.953.03") instead of being auto-coerced into
# numbers by Excel, which would silently drop meaningful trailing zeros / dots.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.953.03'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '3.453.44'
$ws.Range("E3").Value = '  -1.14%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '591.74'
$ws.Range("E5").Value = '  -1.22%  '

$ws.Range("D6").Value = '179.57'
$ws.Range("E6").Value = '  +2.62%  '

$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  +4.02%  '

$ws.Range("D9").Value = '3.450.03'
$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("E10").Value = '  +6.49%  '

$ws.Range("E11").Value = '  -2.47%  '

$ws.Range("D12").Value = '0.429'
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("D13").Value = '4.049.66'
$ws.Range("E13").Value = '  -1.16%  '

$ws.Range("D14").Value = '31.89'
$ws.Range("E14").Value = '  +2.05%  '

$ws.Range("E15").Value = '  -0.70%  '

$ws.Range("D16").Value = '66.938.81'
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("D18").Value = '3.454.81'
$ws.Range("E18").Value = '  -0.95%  '

$ws.Range("D19").Value = '6.19'
$ws.Range("E19").Value = '  -1.38%  '

$ws.Range("D20").Value = '14.15'
$ws.Range("E20").Value = '  -2.52%  '

$ws.Range("D21").Value = '391.03'
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D22").Value = '7.90'
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.27%  '

$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '5.76'
$ws.Range("E24").Value = '  +1.03%  '

$ws.Range("D25").Value = '0.537'
$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("D26").Value = '71.58'
$ws.Range("E26").Value = '  -2.34%  '

$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("E28").Value = '  +1.28%  '

$ws.Range("E29").Value = '  -2.99%  '

$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").Value = '6.11'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  -1.49%  '

$ws.Range("E33").Value = '  -0.82%  '

$ws.Range("D34").Value = '23.49'
$ws.Range("E34").Value = '  -0.58%  '

$ws.Range("D35").Value = '7.31'
$ws.Range("E35").Value = '  -0.82%  '

$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").Value = '1.57'
$ws.Range("E37").Value = '  -3.42%  '

$ws.Range("D38").Value = '160.33'
$ws.Range("E38").Value = '  -1.60%  '

$ws.Range("D39").Value = '0.875'
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("D40").Value = '2.79'
$ws.Range("E40").Value = '  +10.27%  '

$ws.Range("E41").Value = '  -3.03%  '

$ws.Range("D42").Value = '6.71'
$ws.Range("E42").Value = '  -4.57%  '

$ws.Range("D43").Value = '4.64'
$ws.Range("E43").Value = '  -0.57%  '

$ws.Range("D44").Value = '26.02'
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("D45").Value = '0.0717'
$ws.Range("E45").Value = '  -1.67%  '

$ws.Range("D46").Value = '2.751.82'
$ws.Range("E46").Value = '  -1.73%  '

$ws.Range("D47").Value = '26.13'
$ws.Range("E47").Value = '  -4.44%  '

$ws.Range("D48").Value = '41.24'
$ws.Range("E48").Value = '  -2.94%  '

$ws.Range("D49").Value = '0.0297'
$ws.Range("E49").Value = '  -0.97%  '

$ws.Range("D50").Value = '323.34'
$ws.Range("E50").Value = '  -4.25%  '

$ws.Range("E51").Value = '  -2.93%  '
